# Scrum Board update: add the previously-missing "Sprint 1" and "Sprint 3"
# rows (the board jumped straight from nothing to "Sprint 2"), and extend the
# wording of the existing "Sprint 2" discussion task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the "Sprint 1" row right before the existing "Sprint 2" block ---
$ws.Rows("3:3").Insert()
$ws.Rows("3:3").ClearFormats()
$ws.Range("A3").Value = "Sprint 1"
$ws.Range("B3").Value = "Each member of the team should read the code in order to undurstant the functionalities already existent"

# --- Insert two "Sprint 3" rows right after the "Sprint 2" block ---
$ws.Rows("7:8").Insert()
$ws.Rows("7:8").ClearFormats()
$ws.Range("A7").Value = "Sprint 3"

# --- Expand the wording of the existing Sprint 2 task ---
$ws.Range("B6").Value = "The members should discuss the features to be implemented across the next sprint"

$ws.Range("B7").Value = "The team should argue which member would participate in the implementation of each functionality"
$ws.Range("A8").Value = "Sprint 3"
$ws.Range("B8").Value = "The team should build the document to be submited with the user stories about the features to be implemented."

# Restore the cursor/selection like the saved workbook shows
[void]$ws.Range("B17").Select()
